# Add data for 2022-04-24
# - Renames the sheet / "through" label from 04-15 to 04-16
# - Updates the 2022 (Total) column values for March, May, and the Total row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-04-16"

# Update the header label in I1 (shared string "2022 (through 04-15)" -> "... 04-16").
$ws.Range("I1").Value = "2022 (through 04-16)"

# Update the 2022 totals column (column I) for the affected rows.
$ws.Range("I3").Value = 141   # March
$ws.Range("I5").Value = 68    # May
$ws.Range("I14").Value = 503  # Total
